# The sheet contains four daily price records (rows 2-5) for
# "Corazón de apio" that need to be re-ordered onto a weekly cadence.
# The underlying data (row 2 <-> row 3 <-> row 5 <-> row 4 <-> row 2, a
# 4-cycle) stays the same, just redistributed across the date rows, so
# we simply overwrite each touched cell with its new value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was 2021-06-18 / Segunda -> now 2021-03-01 / Primera)
$ws.Range("D2").Value = 44267
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 120
$ws.Range("K2").Value = 1500
$ws.Range("L2").Value = 1800
$ws.Range("M2").Value = 1650
$ws.Range("P2").Value = 275

# Row 3 (was 2021-03-01 / Primera -> now 2022-03-01 / Primera, $/paquete)
$ws.Range("D3").Value = 44623
$ws.Range("J3").Value = 300
$ws.Range("K3").Value = 1800
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = 1900
$ws.Range("N3").Value = "`$/paquete"
$ws.Range("P3").Value = 1900
$ws.Range("Q3").Value = 1

# Row 4 (was 2021-06-11 / Segunda -> now 2021-06-18 / Segunda)
$ws.Range("D4").Value = 44377
$ws.Range("J4").Value = 550
$ws.Range("K4").Value = 2000
$ws.Range("L4").Value = 2800
$ws.Range("M4").Value = 2364
$ws.Range("P4").Value = 394

# Row 5 (was 2022-03-01 / Primera / $/paquete -> now 2021-06-11 / Segunda / $/docena de matas)
$ws.Range("D5").Value = 44370
$ws.Range("I5").Value = "Segunda"
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 1200
$ws.Range("M5").Value = 1080
$ws.Range("N5").Value = "`$/docena de matas"
$ws.Range("P5").Value = 180
$ws.Range("Q5").Value = 6
